$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 79: fill in resultado/profit ---
$ws.Cells.Item(79, 7).Value = "Fallo"
$ws.Cells.Item(79, 8).Value = -1

# --- Append new rows 82-86 ---
$newRows = @(
    @{ Row = 82; A = 14762031; B = "2025-10-04"; C = "Luciano Darderi"; D = "Yunchaokete Bu"; E = "Gana Luciano Darderi"; F = 2 },
    @{ Row = 83; A = 14762048; B = "2025-10-04"; C = "Karen Khachanov"; D = "Juncheng Shang"; E = "Gana Karen Khachanov"; F = 1.5 },
    @{ Row = 84; A = 14785809; B = "2025-10-03"; C = "Pablo Carreño Busta"; D = "Martin Landaluce"; E = "Gana Martin Landaluce"; F = 2 },
    @{ Row = 85; A = 14816745; B = "2025-10-03"; C = "Pedro Sakamoto"; D = "Alan Fernando Rubio Fierros"; E = "Gana Alan Fernando Rubio Fierros"; F = 3 },
    @{ Row = 86; A = 14816687; B = "2025-10-03"; C = "Gabriele Piraino"; D = "Michele Ribecai"; E = "Gana Michele Ribecai"; F = 2.2 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    # Force column B to text so date-like strings aren't converted to date serials
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
}
